$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The three data rows (2-4) shift "down" one slot: a brand-new contact
# (beckjstevens) is written into row 2, the former row-2 contact
# (beatrixpottersociety) moves into row 3, the former row-3 contact
# (medicalmedium) moves into row 4, and the former row-4 contact
# (petermlawrence) drops out entirely. Row/column structure (table range,
# dimension, styles) stays put - only cell contents change, plus one new
# cell L2 ("blog").
# ---------------------------------------------------------------------------

# Row 2 -> beckjstevens (new contact)
$ws.Range("D2").Value = "beckjstevens"
$ws.Range("F2").Value = "Rebecca Stevens ~  Fashion tips + Reels"
$ws.Range("I2").Value = "http://www.jasonstevensproductions.com/"
$ws.Range("J2").Value = "https://www.instagram.com/beckjstevens/"
$ws.Range("K2").Value = 27761
$ws.Range("L2").Value = "blog"
$ws.Range("N2").Value = "UNKNOWN"
$ws.Range("O2").Value = "Stylist /Blogger/ TV Presenter`n🎥 @finallymemovie `n👨 @bigjasonstevens `n👦🏼 @littlepresleystevens `nbeckjstevensmail@gmail.com`n1 Cor 2:9"

# Row 3 -> beatrixpottersociety (was row 2)
$ws.Range("D3").Value = "beatrixpottersociety"
$ws.Range("F3").Value = "The Beatrix Potter Society"
$ws.Range("I3").Value = "https://linktr.ee/thebeatrixpottersociety"
$ws.Range("J3").Value = "https://www.instagram.com/beatrixpottersociety/"
$ws.Range("K3").Value = 14453
$ws.Range("N3").Value = "UNKNOWN"
$ws.Range("O3").Value = "The Beatrix Potter Society is a registered charity and exists to promote the study and appreciation of the life and works of Beatrix Potter. Join us!"

# Row 4 -> medicalmedium (was row 3)
$ws.Range("D4").Value = "medicalmedium"
$ws.Range("F4").Value = "Medical Medium®"
$ws.Range("I4").Value = "http://linktr.ee/medicalmedium"
$ws.Range("J4").Value = "https://www.instagram.com/medicalmedium/"
$ws.Range("K4").Value = 4081267
$ws.Range("N4").Value = "UNKNOWN"
$ws.Range("O4").Value = "📚#1 NY Times Bestselling Author`n🎤Medical Medium Podcast `n🙏🏼Helping people overcome illness `n💚Originator of Global Celery Juice Movement"

# Multi-line Biography text makes the engine auto-grow the row height; snap
# rows 2 and 4 (the ones whose Biography text changed) back with AutoFit so
# no stray ht=/customHeight= survives on the <row> elements.
$ws.Rows(2).AutoFit()
$ws.Rows(4).AutoFit()

# ---------------------------------------------------------------------------
# Hyperlinks: rebuild from scratch in the right order so they line up with
# I2/J2/I3/J3/I4/J4 and their (new) target addresses. Deleting any one
# range's Hyperlinks clears the whole sheet's collection in this engine, so
# do that once up front, then re-add all six.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("I2"), "http://www.jasonstevensproductions.com/")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://www.instagram.com/beckjstevens/")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://linktr.ee/thebeatrixpottersociety")
$ws.Hyperlinks.Add($ws.Range("J3"), "https://www.instagram.com/beatrixpottersociety/")
$ws.Hyperlinks.Add($ws.Range("I4"), "http://linktr.ee/medicalmedium")
$ws.Hyperlinks.Add($ws.Range("J4"), "https://www.instagram.com/medicalmedium/")

# ---------------------------------------------------------------------------
# Hyperlinks.Add() re-stamps the "Hyperlink" cell style via a freshly minted
# cellXf (applyFont="1") instead of reusing the workbook's existing
# Hyperlink-style index, which would leave a near-duplicate style behind.
# Copy the known-good style (still intact on, e.g., K1) back onto the six
# linked cells so they keep pointing at the original style index.
# ---------------------------------------------------------------------------
$holder = $ws.Range("Q1")
$holder.Value = "x"
$holder.Style = "Hyperlink"
$holder.Copy()
$ws.Range("I2:J4").PasteSpecial(-4122)
$holder.Clear()
$excel.CutCopyMode = 0
